$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.874436
$ws.Cells.Item(2, 8).Value = 5.623308
$ws.Cells.Item(2, 9).Value = 0.1442186763702422
$ws.Cells.Item(2, 10).Value = 0.1442186763702422
$ws.Cells.Item(2, 13).Value = 38.745275
$ws.Cells.Item(2, 14).Value = 116.235825
$ws.Cells.Item(2, 15).Value = 0.3160319337595895
$ws.Cells.Item(2, 16).Value = 0.3160319337595895
$ws.Cells.Item(2, 17).Value = 72.62553828989999
$ws.Cells.Item(2, 18).Value = 653.6298446090999
$ws.Cells.Item(2, 19).Value = 0.04557770717753606
$ws.Cells.Item(2, 20).Value = 0.04557770717753604
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.874436
$ws.Cells.Item(3, 8).Value = 5.623308
$ws.Cells.Item(3, 9).Value = 0.1442186763702422
$ws.Cells.Item(3, 10).Value = 0.1442186763702422
$ws.Cells.Item(3, 15).Value = 0.4383510712400526
$ws.Cells.Item(3, 16).Value = 0.4383510712400526
$ws.Cells.Item(3, 17).Value = 100.7350179143
$ws.Cells.Item(3, 18).Value = 906.6151612287
$ws.Cells.Item(3, 19).Value = 0.06321841127971813
$ws.Cells.Item(3, 20).Value = 0.06321841127971811
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.874436
$ws.Cells.Item(4, 8).Value = 5.623308
$ws.Cells.Item(4, 9).Value = 0.1442186763702422
$ws.Cells.Item(4, 10).Value = 0.1442186763702422
$ws.Cells.Item(4, 13).Value = 17.38482166666667
$ws.Cells.Item(4, 14).Value = 52.154465
$ws.Cells.Item(4, 15).Value = 0.1418020341675798
$ws.Cells.Item(4, 16).Value = 0.1418020341675798
$ws.Cells.Item(4, 17).Value = 32.58673558558
$ws.Cells.Item(4, 18).Value = 293.28062027022
$ws.Cells.Item(4, 19).Value = 0.02045050167425622
$ws.Cells.Item(4, 20).Value = 0.02045050167425622
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.874436
$ws.Cells.Item(5, 8).Value = 5.623308
$ws.Cells.Item(5, 9).Value = 0.1442186763702422
$ws.Cells.Item(5, 10).Value = 0.1442186763702422
$ws.Cells.Item(5, 13).Value = 12.72763533333333
$ws.Cells.Item(5, 14).Value = 38.182906
$ws.Cells.Item(5, 15).Value = 0.103814960832778
$ws.Cells.Item(5, 16).Value = 0.103814960832778
$ws.Cells.Item(5, 17).Value = 23.857137863672
$ws.Cells.Item(5, 18).Value = 214.714240773048
$ws.Cells.Item(5, 19).Value = 0.01497205623873177
$ws.Cells.Item(5, 20).Value = 0.01497205623873177
$ws.Cells.Item(6, 9).Value = 0.2460517715407892
$ws.Cells.Item(6, 10).Value = 0.2460517715407892
$ws.Cells.Item(6, 13).Value = 38.745275
$ws.Cells.Item(6, 14).Value = 116.235825
$ws.Cells.Item(6, 15).Value = 0.3160319337595895
$ws.Cells.Item(6, 16).Value = 0.3160319337595895
$ws.Cells.Item(6, 17).Value = 123.906575799225
$ws.Cells.Item(6, 18).Value = 1115.159182193025
$ws.Cells.Item(6, 19).Value = 0.07776021716500835
$ws.Cells.Item(6, 20).Value = 0.07776021716500833
$ws.Cells.Item(7, 9).Value = 0.2460517715407892
$ws.Cells.Item(7, 10).Value = 0.2460517715407892
$ws.Cells.Item(7, 15).Value = 0.4383510712400526
$ws.Cells.Item(7, 16).Value = 0.4383510712400526
$ws.Cells.Item(7, 19).Value = 0.1078570576354176
$ws.Cells.Item(7, 20).Value = 0.1078570576354176
$ws.Cells.Item(8, 9).Value = 0.2460517715407892
$ws.Cells.Item(8, 10).Value = 0.2460517715407892
$ws.Cells.Item(8, 13).Value = 17.38482166666667
$ws.Cells.Item(8, 14).Value = 52.154465
$ws.Cells.Item(8, 15).Value = 0.1418020341675798
$ws.Cells.Item(8, 16).Value = 0.1418020341675798
$ws.Cells.Item(8, 17).Value = 55.596294608745
$ws.Cells.Item(8, 18).Value = 500.366651478705
$ws.Cells.Item(8, 19).Value = 0.03489064171502054
$ws.Cells.Item(8, 20).Value = 0.03489064171502053
$ws.Cells.Item(9, 9).Value = 0.2460517715407892
$ws.Cells.Item(9, 10).Value = 0.2460517715407892
$ws.Cells.Item(9, 13).Value = 12.72763533333333
$ws.Cells.Item(9, 14).Value = 38.182906
$ws.Cells.Item(9, 15).Value = 0.103814960832778
$ws.Cells.Item(9, 16).Value = 0.103814960832778
$ws.Cells.Item(9, 17).Value = 40.702710515658
$ws.Cells.Item(9, 18).Value = 366.324394640922
$ws.Cells.Item(9, 19).Value = 0.02554385502534266
$ws.Cells.Item(9, 20).Value = 0.02554385502534266
$ws.Cells.Item(10, 7).Value = 6.825289333333334
$ws.Cells.Item(10, 8).Value = 20.475868
$ws.Cells.Item(10, 9).Value = 0.5251361975000832
$ws.Cells.Item(10, 10).Value = 0.5251361975000832
$ws.Cells.Item(10, 13).Value = 38.745275
$ws.Cells.Item(10, 14).Value = 116.235825
$ws.Cells.Item(10, 15).Value = 0.3160319337595895
$ws.Cells.Item(10, 16).Value = 0.3160319337595895
$ws.Cells.Item(10, 17).Value = 264.4477121745667
$ws.Cells.Item(10, 18).Value = 2380.0294095711
$ws.Cells.Item(10, 19).Value = 0.165959807983109
$ws.Cells.Item(10, 20).Value = 0.165959807983109
$ws.Cells.Item(11, 7).Value = 6.825289333333334
$ws.Cells.Item(11, 8).Value = 20.475868
$ws.Cells.Item(11, 9).Value = 0.5251361975000832
$ws.Cells.Item(11, 10).Value = 0.5251361975000832
$ws.Cells.Item(11, 15).Value = 0.4383510712400526
$ws.Cells.Item(11, 16).Value = 0.4383510712400526
$ws.Cells.Item(11, 17).Value = 366.8013435847445
$ws.Cells.Item(11, 18).Value = 3301.2120922627
$ws.Cells.Item(11, 19).Value = 0.2301940147210893
$ws.Cells.Item(11, 20).Value = 0.2301940147210893
$ws.Cells.Item(12, 7).Value = 6.825289333333334
$ws.Cells.Item(12, 8).Value = 20.475868
$ws.Cells.Item(12, 9).Value = 0.5251361975000832
$ws.Cells.Item(12, 10).Value = 0.5251361975000832
$ws.Cells.Item(12, 13).Value = 17.38482166666667
$ws.Cells.Item(12, 14).Value = 52.154465
$ws.Cells.Item(12, 15).Value = 0.1418020341675798
$ws.Cells.Item(12, 16).Value = 0.1418020341675798
$ws.Cells.Item(12, 17).Value = 118.6564378834022
$ws.Cells.Item(12, 18).Value = 1067.90794095062
$ws.Cells.Item(12, 19).Value = 0.07446538102053975
$ws.Cells.Item(12, 20).Value = 0.07446538102053975
$ws.Cells.Item(13, 7).Value = 6.825289333333334
$ws.Cells.Item(13, 8).Value = 20.475868
$ws.Cells.Item(13, 9).Value = 0.5251361975000832
$ws.Cells.Item(13, 10).Value = 0.5251361975000832
$ws.Cells.Item(13, 13).Value = 12.72763533333333
$ws.Cells.Item(13, 14).Value = 38.182906
$ws.Cells.Item(13, 15).Value = 0.103814960832778
$ws.Cells.Item(13, 16).Value = 0.103814960832778
$ws.Cells.Item(13, 17).Value = 86.86979367915646
$ws.Cells.Item(13, 18).Value = 781.8281431124082
$ws.Cells.Item(13, 19).Value = 0.05451699377534509
$ws.Cells.Item(13, 20).Value = 0.05451699377534509
$ws.Cells.Item(14, 7).Value = 1.099475
$ws.Cells.Item(14, 8).Value = 3.298425
$ws.Cells.Item(14, 9).Value = 0.08459335458888541
$ws.Cells.Item(14, 10).Value = 0.08459335458888539
$ws.Cells.Item(14, 13).Value = 38.745275
$ws.Cells.Item(14, 14).Value = 116.235825
$ws.Cells.Item(14, 15).Value = 0.3160319337595895
$ws.Cells.Item(14, 16).Value = 0.3160319337595895
$ws.Cells.Item(14, 17).Value = 42.599461230625
$ws.Cells.Item(14, 18).Value = 383.3951510756249
$ws.Cells.Item(14, 19).Value = 0.0267342014339361
$ws.Cells.Item(14, 20).Value = 0.02673420143393609
$ws.Cells.Item(15, 7).Value = 1.099475
$ws.Cells.Item(15, 8).Value = 3.298425
$ws.Cells.Item(15, 9).Value = 0.08459335458888541
$ws.Cells.Item(15, 10).Value = 0.08459335458888539
$ws.Cells.Item(15, 15).Value = 0.4383510712400526
$ws.Cells.Item(15, 16).Value = 0.4383510712400526
$ws.Cells.Item(15, 17).Value = 59.08744487479167
$ws.Cells.Item(15, 18).Value = 531.7870038731249
$ws.Cells.Item(15, 19).Value = 0.03708158760382754
$ws.Cells.Item(15, 20).Value = 0.03708158760382753
$ws.Cells.Item(16, 7).Value = 1.099475
$ws.Cells.Item(16, 8).Value = 3.298425
$ws.Cells.Item(16, 9).Value = 0.08459335458888541
$ws.Cells.Item(16, 10).Value = 0.08459335458888539
$ws.Cells.Item(16, 13).Value = 17.38482166666667
$ws.Cells.Item(16, 14).Value = 52.154465
$ws.Cells.Item(16, 15).Value = 0.1418020341675798
$ws.Cells.Item(16, 16).Value = 0.1418020341675798
$ws.Cells.Item(16, 17).Value = 19.11417680195833
$ws.Cells.Item(16, 18).Value = 172.027591217625
$ws.Cells.Item(16, 19).Value = 0.01199550975776333
$ws.Cells.Item(16, 20).Value = 0.01199550975776332
$ws.Cells.Item(17, 7).Value = 1.099475
$ws.Cells.Item(17, 8).Value = 3.298425
$ws.Cells.Item(17, 9).Value = 0.08459335458888541
$ws.Cells.Item(17, 10).Value = 0.08459335458888539
$ws.Cells.Item(17, 13).Value = 12.72763533333333
$ws.Cells.Item(17, 14).Value = 38.182906
$ws.Cells.Item(17, 15).Value = 0.103814960832778
$ws.Cells.Item(17, 16).Value = 0.103814960832778
$ws.Cells.Item(17, 17).Value = 13.99371685811667
$ws.Cells.Item(17, 18).Value = 125.94345172305
$ws.Cells.Item(17, 19).Value = 0.008782055793358437
$ws.Cells.Item(17, 20).Value = 0.008782055793358435

Write-Output "applied changes"